$d = $word.ActiveDocument

$pairs = @(
    @{ old = "2025-05-03 Saturday"; new = "2025-05-04 Sunday" },
    @{ old = "717×2="; new = "433×2=" },
    @{ old = "285×9="; new = "904×9=" },
    @{ old = "845×7="; new = "154×9=" },
    @{ old = "643×5="; new = "875×5=" },
    @{ old = "888×6="; new = "885×8=" },
    @{ old = "889×2="; new = "636×2=" },
    @{ old = "410×3="; new = "391×5=" },
    @{ old = "962×7="; new = "499×7=" },
    @{ old = "734×5="; new = "926×2=" },
    @{ old = "589×2="; new = "533×3=" },
    @{ old = "132×8="; new = "853×3=" },
    @{ old = "450×6="; new = "164×5=" },
    @{ old = "305×2="; new = "275×3=" },
    @{ old = "657×7="; new = "440×5=" },
    @{ old = "977×6="; new = "684×7=" },
    @{ old = "979×9="; new = "611×3=" },
    @{ old = "251×3="; new = "310×2=" },
    @{ old = "884×6="; new = "367×3=" },
    @{ old = "202×6="; new = "831×6=" },
    @{ old = "275×5="; new = "818×8=" },
    @{ old = "466×3="; new = "434×7=" },
    @{ old = "839×9="; new = "741×6=" },
    @{ old = "923×6="; new = "531×4=" },
    @{ old = "247×7="; new = "328×6=" },
    @{ old = "831×7="; new = "388×8=" }
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}
